$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet from "My Series" to "Data"
$ws.Name = "Data"

# 2. Widen the custom number format for numFmtId 164 (used by A2:A11)
#    from "0.000" to "###0.000"
$ws.Range("A2:A11").NumberFormat = "###0.000"

# 3. Replace the cell-A1 comment text (opaque CEIC add-in metadata blob)
#    with the updated payload
$cmt = $ws.Range("A1").Comment
[void]$cmt.Text("GRoAAB+LCAAAAAAAAAOlGdtu48b1Vwg9tUAlktLa8RqzDHSzI1SyDEmu47wUI3IsTU2RCjm0rLcUaJEiTVEUxaZIr+hTigLdLtoESHd7+Zdg7d0+9Rd65sLhRXK22i6MNefc5syZcx2jt28WvnFNopiGwaOKXbMqBgnc0KPB7FElYZdVe7/ytoO6Ny7xT3GEF4QBsQFcQXx4E9NHlTljy0PTXK1WtVWjFkYzs25ZtvnuoD9252SBqzSIGQ5cUtFc3uu5Kg5qe4sBYdjDDEvOR5XeuFdrE+p2ADbAAZ6RqNZKYhqQOO4GjDJKYs4ZEcxIuzP4jjyYU6/t12xkbsAzylZCfU/SFSglXNHBtmRCF8SpW/ZB1TqoNuyJbR/W9w8bjdrewf57KaMmRH0cszGJrqkrAGOGF0vBbh00bNuu7zcayNxKBLIyAzho6Hsjck1j4rWJ78c7WcRUF9h0GZx6N2NayMzxKkFvrsJxhJfzCWU+2U2N0aBlLAKlSybEQUdhRFyw3xupdEJWw0iZdbLsA3YypxFbd/B6Z1lnMYmGS26k3Vgd1AkD1vRJxM6WcNfEA1cAhMOihCDzHmTG1KGxC980SIjnXGI/zjMVkOg8jK7iJXbJCcSxyWWsAj/EHjgcozGjbpwJ2MCg0yhcgkjYvRX63hGIVSpuQWjRvQBszPdtheFVWXoRicS1iguGS11glpJvwNF4Hq6Ggb8eJ9PYjeiUeJ1WSr0Vh3hEKu52ErNwAVpkICRhOcga/kEElsGoQ1y6wP6pD1aMHYjdIgA1ExZeUtYO/WQRaHOWoOgcTjQhN/qEeo2GcLsBN3oY9IKUXpp5K6rIMApXes9NhDBCDtyM3dTJNhFl4g7A0uvbxIgb4ac8oj4UiPxd5KBFrxjPCWFbXUJiEM+FR7zkOK31SbKYQnxNIciuxa4xMjM8Aj8FZwe9HAvKSFX8TCzrUPyAHhqNuoF3P12KRLBdbi/HBlwJhOBMfsvHwRVAzymbnzTTs2zBIGmBe+k3cQhCd+njtQBrK+VhqBe4fuIRmRF6waVwUa6bor4XjzZAfQhyB+FgPVkveWYw76HokEuc+FCsGOSUWZYuSmDUjK/KNHkQOov89MYd3grE0Au43qLmQrbk9a7mhgsOMKEEn4+RmafnKd8l3WDWx8EsgaSq7ViGa3/j+WAS4SDmx9E5tOR624lQei8ytzvysoaJcCd5WSFgkVmiQxOyWIYR9gdgGHqUBKJ+qgIBFh5gNlcriGWfuKmRzYxVcxU1SxV/HZkISnkMnvBVrJeAgoifRXYducDNgIgfcxB6UPixT6dRwcu24uDGsnKY+hs/3Y6lMb0EaD0h3XybrHk/ki0UnDutY6cIsUJtUMkZjx4c1PesRh1SOV8jceQRwb7RhaaUEaMXXJOYLYDt0BiRmHrwRbF/aLxDpoRC1AsbqbjbmTvPh44i8n4CzfVaqNIEfYuQIgEk2BkNsL9JqDEZg3NBcOSvc4TyqP3QBbq7H/3r9tdPXzz77O6jx6+++MF//v7LF//46e2TH8LH3V/+evvxL+QxJTGa4KlPhEKT1sGB1XgAjqZBSOUHKPle4jIBu7gQLYBeI9XKikW722sf91sioWhgys67i4QHQx+vwyRbjuUhxEbiSs3UEySJM0kTlFoXsB3C676MQt61XpMidR5/H6O0xcvnn718/qd7uZXBsuJiP3y4V7Xrr6090P/bG3S69vBpQDZ76azxoGrtVev1HHGJBo1gZILeTtup5zkN23po1RuWrZO5px15G1EZpSRN8Mws8UlQO0wCFq21C+TXKVI4/gRCRKNlKOQWykW/+PGrPz8uUCnrKkhRCigXJpH0JjNdCNEno4kxHp6N2l1j0h1zP8lwOTop/GuI1e46ngpOFQQJ9r8FgzIfk40KtIAVI7w0CHbnxhoiMReHBWfbBpUbvaHIspbHUZgs5Y3kGDLoFkqdTbZybMk1AifsuZF0MtQWcqnr7d8+38agDqLc7CygTE98eRgqYCQoh1dR++k/X3z54Ytnz+6e/uz2y+8XJKh99CAEfg7RlF9qt4eUp+pNCYLOx8KYV9Z3c/VFAXnneBrSgMWOvS+aRrVCwGpzaeI36i2g5AnBwl4AL0HQOzju3jAV2M4JMosA0HOJodqGWbOtATKHZ3b9929+e/erz+8+efrqwz/efvSH248/efn8d6+e/F5G3d3jp3c/eaKyfLkQCF14Cy+7QEMMhK7Bo9Hgtdv46oOfG0HIDOg5jERkpK8++DQnjCsqupNMMvR0WpGiChukeWbOZ+RU0ToU+DSLbADavIQ1NIUqYuGSutkm71W5KB53AvGN3qSaxMQIoZ36JpykSJwx/698ikWW1NO3rLpdV1ipDT/CFMc50x/74RSajBQhJq4SSYHr6xkyWrHfcX/YavYzEqnEMPJgQLP42Mg/UNpT8pLSi9OVniwyCGCh8XMTnw/JG2SbKC05l8ZMNW1eNj2e/rbPZgUKGNijSDZEgXq9HCdL6IbTN4n78eKhJtcAn8heNd8SZ+tep4iHdQ4LhbCI5gCBF6lJoWSa6sV8sJXt7Ak3TbYEXOFxB8yhHihlp3UNfWVk8rzTjaIw2pp8MkxKNoBOGjKKmVlc04g7lV23l91VCkgTnv6Qo586YdghPmG7vd6ZGfcgvH5jXrj7XVl78dD3lDF3Gz20WTIB+SdM7ij/7wumdLZmFEFjxV88dn5yTCfXEQy8O2ojjyIY+QgIu6s3wSMaxexdngnUl4RcaMiF7FAFifwQ6wvxSCK/1CGVdLOgZhq6TD4kh36fLuiOY6GVxndRCNhyuZQtXG83T+Gl5YTcQIOZkwBJcfo9KBt85NlNmnRYyKWanz/WxHQ2Z7sq9tYUE49Mrao7JfXqA886qD4kpFG1bfgfu/W6Ze3xlx4lHDIHJasdNzHTC8v+wOP8FysPGvUZGgAA")
